$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Result Summary" cell - insert "RESULTS_SUMMARY_IN " before the
# "See comment below regarding the test limitations..." run.
# ---------------------------------------------------------------------------

# Locate the unique anchor text right before the insertion point so the
# positions below are computed from the live document rather than hard-coded
# offsets.
$anchor = $d.Content
$anchor.Find.Execute("Failed assay due to suboptimal DNA quantity/quality")
$failedEnd = $anchor.End

# Layout immediately after the anchor: "." then " " then "See comment..."
$dotStart = $failedEnd
$dotEnd = $dotStart + 1
$spaceStart = $dotEnd
$spaceEnd = $spaceStart + 1
$seeStart = $spaceEnd

$newText = "RESULTS_SUMMARY_IN "
$newTextNoSpaceLen = ("RESULTS_SUMMARY_IN").Length

# Insert the new text immediately before "See comment..." (this naturally
# merges with whatever run precedes the insertion point).
$insPoint = $d.Range($seeStart, $seeStart)
$insPoint.InsertBefore($newText)

$newTextEnd = $seeStart + $newText.Length
$newTextMid = $seeStart + $newTextNoSpaceLen

# Split the freshly inserted text into its own two runs (plain, non-bold,
# Arial complex-script font) - format the trailing segment first so a run
# boundary is created between the two pieces.
$segSpace = $d.Range($newTextMid, $newTextEnd)
$segSpace.Font.Bold = 0
$segSpace.Font.NameBi = "Arial"

$segWord = $d.Range($seeStart, $newTextMid)
$segWord.Font.Bold = 0
$segWord.Font.NameBi = "Arial"

# Restore the original "." / " " run split that the insertion step merged
# together, by re-asserting (toggling) their formatting so a run boundary is
# recreated between them and ahead of the new text.
$segSpaceOrig = $d.Range($spaceStart, $spaceEnd)
$segSpaceOrig.Font.Bold = 0
$segSpaceOrig.Font.Bold = 1

$segDot = $d.Range($dotStart, $dotEnd)
$segDot.Font.Bold = 0
$segDot.Font.Bold = 1

# ---------------------------------------------------------------------------
# Change 2: update the cached SAVEDATE field result from 15-Mar-2024 to
# 6-Sep-2024.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("15-Mar-2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "6-Sep-2024", 2)
